$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 <-> Row 4 swap of the Protocol Names / Folder Index / Colors columns (G, H, I)
$ws.Range("G2").Value = "4 degree height"
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = "black"

$ws.Range("H3").Value = 2

$ws.Range("G4").Value = "16 degree height"
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = "red"

# New label in C4 (Axes Names column)
$ws.Range("C4").Value = "Absolute Angle of Roll (|°|)"

# Update the active selection to I5
$ws.Range("I5").Select()
